# edit.ps1
# Applies the three changes described by the commit "Updated Categoria with
# PUT and some repairs":
#
#  1. Paragraph "Camada de domínio: classes de domínio" gets a new trailing
#     run ", nada mais são que as Entidades." and paragraph "Camada de
#     serviço: ..." gets a new trailing run "..algo que não foi implementado
#     nas classes de domínio."
#  2. The paragraph "Em CategoriaResources" loses the <w:lang w:val="en-US"/>
#     from both its paragraph-mark rPr and its run rPr.
#  3. The three runs making up "Após finalizar a Classe, abrir o Postman ..."
#     are merged into a single run, and two new empty paragraphs are added
#     right after it - the first with a bottom paragraph border, the second
#     plain (both carrying sz/szCs = 20).

$d = $word.ActiveDocument

function Find-ParagraphIndex {
    param([string]$NeedleText)
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like $NeedleText) {
            return $i
        }
    }
    return -1
}

function Add-TrailingRun {
    # Appends $Text as a brand-new <w:r> (not merged into the last existing
    # run) at the end of the paragraph at 1-based index $ParaIndex.
    param([int]$ParaIndex, [string]$Text)

    $p = $d.Paragraphs.Item($ParaIndex)
    $insertPos = $p.Range.End - 1
    $mark = $d.Range($insertPos, $insertPos)
    $null = $mark.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($ParaIndex + 1)
    $null = $newPara.Range.InsertAfter($Text)

    # Remove the paragraph break we introduced so the new text re-joins the
    # original paragraph as a second, separate run.
    $breakRange = $d.Range($insertPos, $insertPos + 1)
    $null = $breakRange.Delete()
}

# ---------------------------------------------------------------------
# Change 1: add the two trailing runs.
# ---------------------------------------------------------------------

$idxDominio = Find-ParagraphIndex "Camada de domínio: classes de domínio*"
Add-TrailingRun $idxDominio ", nada mais são que as Entidades."

$idxServico = Find-ParagraphIndex "Camada de serviço: oferece consultas*"
Add-TrailingRun $idxServico "..algo que não foi implementado nas classes de domínio."

# ---------------------------------------------------------------------
# Change 2: drop the en-US language tag on "Em CategoriaResources".
# ---------------------------------------------------------------------

$idxCategoriaResources = Find-ParagraphIndex "Em CategoriaResources*"
$p = $d.Paragraphs.Item($idxCategoriaResources)
$full = $p.Range
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Em CategoriaResources</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $full.InsertXML($xml)

# ---------------------------------------------------------------------
# Change 3: merge the Postman/localhost runs and append the two new
# paragraphs (one with a bottom border, one plain).
# ---------------------------------------------------------------------

$idxPostman = Find-ParagraphIndex "*finalizar a Classe, abrir o Postman*"

# Keep a guard paragraph past the end of the document at every step so the
# range we are about to rewrite with InsertXML is never the document's
# very-last paragraph (that special case inserts one extra empty paragraph).
$guardAnchor = $d.Paragraphs.Item($idxPostman).Range.End
$null = $d.Range($guardAnchor, $guardAnchor).InsertParagraphAfter()

$pPostman = $d.Paragraphs.Item($idxPostman)
$fullPostman = $pPostman.Range
$xmlMerged = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Após </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>finalizar a Classe, abrir o Postman e dar um Post : localhost:8080/categorias no raw e escolhido o JSON, depois colocado URI com o endpoint e dado get, ok deu certo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $fullPostman.InsertXML($xmlMerged)

# Placeholder paragraph for the bordered paragraph; add another guard so the
# placeholder is never the last paragraph when we rewrite it below.
$pPostman2 = $d.Paragraphs.Item($idxPostman)
$afterPostman = $pPostman2.Range.End
$null = $d.Range($afterPostman, $afterPostman).InsertParagraphAfter()

$pBorderPlaceholder = $d.Paragraphs.Item($idxPostman + 1)
$afterBorderPlaceholder = $pBorderPlaceholder.Range.End
$null = $d.Range($afterBorderPlaceholder, $afterBorderPlaceholder).InsertParagraphAfter()

# Fill in the bordered paragraph.
$pBorder = $d.Paragraphs.Item($idxPostman + 1)
$fullBorder = $pBorder.Range
$xmlBorder = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $fullBorder.InsertXML($xmlBorder)

# Fill in the plain trailing paragraph (still guarded by the extra empty
# paragraph appended after it).
$pPlain = $d.Paragraphs.Item($idxPostman + 2)
$fullPlain = $pPlain.Range
$xmlPlain = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $fullPlain.InsertXML($xmlPlain)

# Remove the final leftover guard paragraph (always empty, always the very
# last paragraph of the document at this point).
$lastIndex = $d.Paragraphs.Count
$pGuard = $d.Paragraphs.Item($lastIndex)
$guardDelete = $d.Range($pGuard.Range.Start - 1, $pGuard.Range.End)
$null = $guardDelete.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
